$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.520.89'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '3.012.15'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.16'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.97%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -2.30%  '
$ws.Range('D9').Value = '3.010.36'
$ws.Range('E9').Value = '  -1.30%  '
$ws.Range('E10').Value = '  -3.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.80'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('E12').Value = '  +3.22%  '
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.63'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.85%  '
$ws.Range('E15').Value = '  +2.23%  '
$ws.Range('D16').Value = '3.510.97'
$ws.Range('E16').Value = '  -1.18%  '
$ws.Range('E17').Value = '  -1.16%  '
$ws.Range('D18').Value = '62.480.94'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').Value = '3.009.87'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '460.14'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.691'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.32%  '
$ws.Range('E23').Value = '  -1.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.82'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.46%  '
$ws.Range('E25').Value = '  -8.58%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.69%  '
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.90'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.01'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.76%  '
$ws.Range('E32').Value = '  -4.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.80'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.110'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').Value = '0.0₃0815'
$ws.Range('E35').Value = '  -0.65%  '
$ws.Range('E36').Value = '  -3.26%  '
$ws.Range('E37').Value = '  -2.46%  '
$ws.Range('E38').Value = '  -4.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.21'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '50.40'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.124'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.28%  '
$ws.Range('E42').Value = '  -10.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '392.45'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -10.01%  '
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.269'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.73%  '
$ws.Range('D46').Value = '2.741.88'
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '37.51'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.77'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('E51').Value = '  -0.66%  '
